$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (2-11, columns B:G) down by one row.
# Row 11's original values get overwritten (discarded), and a new row
# of values is written into row 2.
for ($r = 10; $r -ge 2; $r--) {
    $src = $r
    $dst = $r + 1
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($dst, $c).Value = $ws.Cells.Item($src, $c).Value2
    }
}

# Write the new values for row 2 (B2:G2)
$ws.Cells.Item(2, 2).Value = 0.1369420173923726
$ws.Cells.Item(2, 3).Value = 0.9020021672123393
$ws.Cells.Item(2, 4).Value = 4.56720535997291
$ws.Cells.Item(2, 5).Value = 2.137102093951739
$ws.Cells.Item(2, 6).Value = 2.15627664282098
$ws.Cells.Item(2, 7).Value = 46
